$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2039.35
$ws.Range("I15").Value = 2039.35
$ws.Range("K15").Value = 6118.049999999999
$ws.Range("M15").Value = -5949.049999999999
$ws.Range("H86").Value = 2493.7144
$ws.Range("I86").Value = 2145
$ws.Range("J86").Value = 3060.375
$ws.Range("K86").Value = 2145
$ws.Range("L86").Value = 3060.375
$ws.Range("M86").Value = -1022
$ws.Range("N86").Value = -5306.375
$ws.Range("H89").Value = 2493.7144
$ws.Range("I89").Value = 2145
$ws.Range("J89").Value = 3060.375
$ws.Range("K89").Value = 10725
$ws.Range("L89").Value = 15301.875
$ws.Range("M89").Value = -5109
$ws.Range("N89").Value = -26533.875
$ws.Range("H98").Value = 11431.846
$ws.Range("I98").Value = 7534.8887
$ws.Range("K98").Value = 7534.8887
$ws.Range("M98").Value = -6036.8887
$ws.Range("H112").Value = 2318.95
$ws.Range("J112").Value = 2923.4
$ws.Range("L112").Value = 8770.200000000001
$ws.Range("N112").Value = -10986.2
$ws.Range("H122").Value = 11431.846
$ws.Range("I122").Value = 7534.8887
$ws.Range("K122").Value = 22604.6661
$ws.Range("M122").Value = -20154.6661
$ws.Range("H127").Value = 2465.8333
$ws.Range("I127").Value = 930
$ws.Range("J127").Value = 2977.7778
$ws.Range("K127").Value = 2790
$ws.Range("L127").Value = 8933.3334
$ws.Range("M127").Value = 2170
$ws.Range("N127").Value = -18853.3334
$ws.Range("H129").Value = 1188.7736
$ws.Range("I129").Value = 595
$ws.Range("J129").Value = 1212.0588
$ws.Range("K129").Value = 1785
$ws.Range("L129").Value = 3636.1764
$ws.Range("M129").Value = 3215
$ws.Range("N129").Value = -13636.1764
$ws.Range("H138").Value = 2073.7942
$ws.Range("I138").Value = 1207.5
$ws.Range("J138").Value = 6116.5
$ws.Range("K138").Value = 3622.5
$ws.Range("L138").Value = 18349.5
$ws.Range("M138").Value = 1517.5
$ws.Range("N138").Value = -28629.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1212.4445
$ws.Range("I74").Value = 1176.5
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1176.5
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -302.5
$ws.Range("N74").Value = -3248
$ws.Range("H77").Value = 1212.4445
$ws.Range("I77").Value = 1176.5
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 5882.5
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -1514.5
$ws.Range("N77").Value = -16236
$ws.Range("H102").Value = 4123.3887
$ws.Range("I102").Value = 3767.7778
$ws.Range("J102").Value = 4479
$ws.Range("K102").Value = 3767.7778
$ws.Range("L102").Value = 4479
$ws.Range("M102").Value = -2145.7778
$ws.Range("N102").Value = -7723
$ws.Range("H122").Value = 2577.5334
$ws.Range("I122").Value = 2045.375
$ws.Range("J122").Value = 3185.7144
$ws.Range("K122").Value = 6136.125
$ws.Range("L122").Value = 9557.143199999999
$ws.Range("M122").Value = -3686.125
$ws.Range("N122").Value = -14457.1432

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 251250
$ws.Range("I20").Value = 501000
$ws.Range("K20").Value = 501000
$ws.Range("M20").Value = -500753
$ws.Range("H99").Value = 2300
$ws.Range("I99").Value = 2115.3845
$ws.Range("J99").Value = 2642.8572
$ws.Range("K99").Value = 2115.3845
$ws.Range("L99").Value = 2642.8572
$ws.Range("M99").Value = -617.3845000000001
$ws.Range("N99").Value = -5638.8572
$ws.Range("H105").Value = 4929179
$ws.Range("I105").Value = 6805665
$ws.Range("J105").Value = 3403.5
$ws.Range("K105").Value = 6805665
$ws.Range("L105").Value = 3403.5
$ws.Range("M105").Value = -6803918
$ws.Range("N105").Value = -6897.5
$ws.Range("H108").Value = 50684
$ws.Range("J108").Value = 50684
$ws.Range("L108").Value = 50684
$ws.Range("N108").Value = -58364

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1964.4445
$ws.Range("I58").Value = 1775.6154
$ws.Range("J58").Value = 2455.4
$ws.Range("K58").Value = 1775.6154
$ws.Range("L58").Value = 2455.4
$ws.Range("M58").Value = -1572.6154
$ws.Range("N58").Value = -2861.4
$ws.Range("H107").Value = 411.5625
$ws.Range("I107").Value = 342.8889
$ws.Range("K107").Value = 342.8889
$ws.Range("M107").Value = 1577.1111
$ws.Range("H136").Value = 1964.4445
$ws.Range("I136").Value = 1775.6154
$ws.Range("J136").Value = 2455.4
$ws.Range("K136").Value = 5326.8462
$ws.Range("L136").Value = 7366.200000000001
$ws.Range("M136").Value = -2776.8462
$ws.Range("N136").Value = -12466.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 575.6667
$ws.Range("I46").Value = 490.8
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1472.4
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -1381.4
$ws.Range("N46").Value = -3182
$ws.Range("H121").Value = 31735.277
$ws.Range("I121").Value = 1138.5714
$ws.Range("J121").Value = 51205.91
$ws.Range("K121").Value = 3415.7142
$ws.Range("L121").Value = 153617.73
$ws.Range("M121").Value = -2105.7142
$ws.Range("N121").Value = -156237.73
$ws.Range("H129").Value = 3125896.5
$ws.Range("J129").Value = 3572367.2
$ws.Range("L129").Value = 10717101.6
$ws.Range("N129").Value = -10727101.6
$ws.Range("H131").Value = 2416.3118
$ws.Range("I131").Value = 484
$ws.Range("J131").Value = 2883.8064
$ws.Range("K131").Value = 1452
$ws.Range("L131").Value = 8651.4192
$ws.Range("M131").Value = 3588
$ws.Range("N131").Value = -18731.4192
$ws.Range("H141").Value = 7022.5
$ws.Range("I141").Value = 7022.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 21067.5
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -15887.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3412.875
$ws.Range("I80").Value = 2987.5
$ws.Range("J80").Value = 3838.25
$ws.Range("K80").Value = 2987.5
$ws.Range("L80").Value = 3838.25
$ws.Range("M80").Value = -1989.5
$ws.Range("N80").Value = -5834.25
$ws.Range("H83").Value = 3412.875
$ws.Range("I83").Value = 2987.5
$ws.Range("J83").Value = 3838.25
$ws.Range("K83").Value = 14937.5
$ws.Range("L83").Value = 19191.25
$ws.Range("M83").Value = -9945.5
$ws.Range("N83").Value = -29175.25
$ws.Range("H97").Value = 36569.934
$ws.Range("I97").Value = 52963.8
$ws.Range("J97").Value = 3782.2
$ws.Range("K97").Value = 52963.8
$ws.Range("L97").Value = 3782.2
$ws.Range("M97").Value = -52467.8
$ws.Range("N97").Value = -4774.2
$ws.Range("H114").Value = 37722
$ws.Range("J114").Value = 37722
$ws.Range("L114").Value = 37722
$ws.Range("N114").Value = -46400

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4636.154
$ws.Range("I40").Value = 4622.5
$ws.Range("J40").Value = 4658
$ws.Range("K40").Value = 4622.5
$ws.Range("L40").Value = 4658
$ws.Range("M40").Value = -4486.5
$ws.Range("N40").Value = -4930
$ws.Range("H93").Value = 2422.6667
$ws.Range("I93").Value = 2265.3333
$ws.Range("K93").Value = 2265.3333
$ws.Range("M93").Value = -1017.3333
$ws.Range("H122").Value = 150001660
$ws.Range("I122").Value = 125002500
$ws.Range("J122").Value = 200000000
$ws.Range("K122").Value = 375007500
$ws.Range("L122").Value = 600000000
$ws.Range("M122").Value = -375005050
$ws.Range("N122").Value = -600004900
$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800
$ws.Range("H136").Value = 3233.8262
$ws.Range("I136").Value = 2377.8
$ws.Range("J136").Value = 3892.3076
$ws.Range("K136").Value = 7133.400000000001
$ws.Range("L136").Value = 11676.9228
$ws.Range("M136").Value = -4583.400000000001
$ws.Range("N136").Value = -16776.9228

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2282.3809
$ws.Range("I126").Value = 1262.4
$ws.Range("J126").Value = 4832.3335
$ws.Range("K126").Value = 3787.2
$ws.Range("L126").Value = 14497.0005
$ws.Range("M126").Value = -1317.2
$ws.Range("N126").Value = -19437.0005
$ws.Range("H136").Value = 2566.3333
$ws.Range("I136").Value = 2207.9167
$ws.Range("K136").Value = 6623.750100000001
$ws.Range("M136").Value = -4073.750100000001
